$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Began analyze Goldbetter1995 PE analysis: append "Species" fit-item rows
# (initial concentrations) below the existing "ReactionParameter" rows.
$rows = @(
    @{ A = "monophosphorylated PER";  B = "0.25"; E = "reactions";  G = "CN=Root,Model=Goldbeter1995_CircClock,Vector=Compartments[CYTOPLASM],Vector=Metabolites[monophosphorylated PER]" },
    @{ A = "total PER";               B = "1.0";  E = "assignment"; G = "CN=Root,Model=Goldbeter1995_CircClock,Vector=Compartments[CYTOPLASM],Vector=Metabolites[total PER]" },
    @{ A = "PER mRNA";                B = "0.1";  E = "reactions";  G = "CN=Root,Model=Goldbeter1995_CircClock,Vector=Compartments[CYTOPLASM],Vector=Metabolites[PER mRNA]" },
    @{ A = "nuclear PER";             B = "0.25"; E = "reactions";  G = "CN=Root,Model=Goldbeter1995_CircClock,Vector=Compartments[NUCLEUS],Vector=Metabolites[nuclear PER]" },
    @{ A = "EmptySet";                B = "0.0";  E = "fixed";      G = "CN=Root,Model=Goldbeter1995_CircClock,Vector=Compartments[default],Vector=Metabolites[EmptySet]" },
    @{ A = "unphosphorylated PER";    B = "0.25"; E = "reactions";  G = "CN=Root,Model=Goldbeter1995_CircClock,Vector=Compartments[CYTOPLASM],Vector=Metabolites[unphosphorylated PER]" },
    @{ A = "biphosphorylated PER";    B = "0.25"; E = "reactions";  G = "CN=Root,Model=Goldbeter1995_CircClock,Vector=Compartments[CYTOPLASM],Vector=Metabolites[biphosphorylated PER]" }
)

$startRow = 20
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $item = $rows[$i]

    # Column A ("Parameter") carries the bold/centred/bordered header style,
    # same as every other row in the table - copy it from an existing cell.
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item(2, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)

    # StartValue / LowerBound / UpperBound are stored as plain text (not
    # numbers) in this template, so force text storage then drop the
    # number-format override again so no extra style is left behind.
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 2).ClearFormats()

    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = "1e-06"
    $ws.Cells.Item($r, 3).ClearFormats()

    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = "1000000"
    $ws.Cells.Item($r, 4).ClearFormats()

    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = "Species"
    $ws.Cells.Item($r, 7).Value = $item.G
}

$excel.CutCopyMode = 0
